# Update "200 run - Ascended" sheet data (column D, rows 3-13) and
# the active sheet/selection to match the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Optimizer | 200 run - Ascended")

$newValues = @{
    3  = 83877
    4  = 355794
    5  = 693264
    6  = 1067934
    7  = 1479804
    8  = 1928874
    9  = 2415144
    10 = 2938614
    11 = 3499284
    12 = 4097154
    13 = 4712149
}

foreach ($row in $newValues.Keys) {
    $ws.Range("D" + $row).Value = $newValues[$row]
}

# Match the saved UI state: "200 run - Ascended" becomes the active sheet,
# with M17 selected (previously "Optimizer Disabled - Descended" was active
# with F18 selected).
$ws.Activate()
$ws.Range("M17").Select()
